$d = $word.ActiveDocument

# Replace figures in the main document body/tables
$d.Content.Find.Execute("41 770 000,00", $false, $false, $false, $false, $false, $true, 1, $false, "44 260 000,00", 2)
$d.Content.Find.Execute("3 759 300,00", $false, $false, $false, $false, $false, $true, 1, $false, "3 983 400,00", 2)

# Replace the year shown in the textbox/shape
$d.Content.Find.Execute("2020", $false, $false, $false, $false, $false, $true, 1, $false, "2021", 2)
